# Update cryptocurrency price/volume figures to the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.413.24"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "'3.765.97"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'614.65"
$ws.Range("D6").Value = "'177.26"
$ws.Range("E6").Value = "  -2.17%  "
$ws.Range("D7").Value = "'3.764.33"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'6.46"
$ws.Range("E11").Value = "  +2.98%  "
$ws.Range("D12").Value = "'0.483"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("D13").Value = "'39.65"
$ws.Range("E13").Value = "  -4.34%  "
$ws.Range("D14").Value = "'0.0000254"
$ws.Range("E14").Value = "  -1.11%  "
$ws.Range("D15").Value = "'4.397.27"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "'3.767.22"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").Value = "'69.458.74"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "'7.53"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").Value = "'508.52"
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("D21").Value = "'16.38"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").Value = "'9.47"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  +0.82%  "
$ws.Range("D24").Value = "'2.47"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").Value = "'86.23"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("D26").Value = "'0.0000142"
$ws.Range("E26").Value = "  +7.63%  "
$ws.Range("D27").Value = "'12.83"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").Value = "'10.51"
$ws.Range("E28").Value = "  -5.26%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "'2.54"
$ws.Range("E30").Value = "  +1.07%  "
$ws.Range("D31").Value = "'3.01"
$ws.Range("E31").Value = "  +2.89%  "
$ws.Range("D32").Value = "'8.09"
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("D33").Value = "'30.93"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").Value = "'1.04"
$ws.Range("E36").Value = "  -1.70%  "
$ws.Range("D37").Value = "'6.11"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("E38").Value = "  +5.95%  "
$ws.Range("D39").Value = "'0.338"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").Value = "'466.53"
$ws.Range("E40").Value = "  +9.88%  "
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").Value = "'3.01"
$ws.Range("E42").Value = "  +8.79%  "
$ws.Range("D43").Value = "'49.88"
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").Value = "'44.22"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("D45").Value = "'8.56"
$ws.Range("E45").Value = "  -1.49%  "
$ws.Range("D46").Value = "'2.945.39"
$ws.Range("E46").Value = "  -2.34%  "
$ws.Range("D47").Value = "'0.0361"
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'27.37"
$ws.Range("E48").Value = "  -0.64%  "
$ws.Range("D49").Value = "'139.29"
$ws.Range("E49").Value = "  +1.86%  "
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("E51").Value = "  -1.09%  "
